# Weekly update: insert one new price record for "Pimiento" (Cuatro cascos
# verde, Region del Maule) ahead of the existing row 319, shifting the
# remaining records down by one row (old row 319..396 -> new row 320..397).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 319; Excel shifts rows 319-396 down to
# 320-397 and extends the sheet's used range/dimension automatically.
$ws.Rows.Item(319).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A319").Value = 5
$ws.Range("B319").Value = "Macroferia Regional de Talca"
$ws.Range("C319").Value = "Maule"
$ws.Range("D319").Value = 44543
$ws.Range("E319").Value = 7
$ws.Range("F319").Value = 100112002
$ws.Range("G319").Value = "Pimiento"
$ws.Range("H319").Value = "Cuatro cascos verde"
$ws.Range("I319").Value = "Primera"
$ws.Range("J319").Value = 300
$ws.Range("K319").Value = 13000
$ws.Range("L319").Value = 13000
$ws.Range("M319").Value = 13000
$ws.Range("N319").Value = "$/caja 15 kilos"
$ws.Range("O319").Value = "Región del Maule"
$ws.Range("P319").Value = 867
$ws.Range("Q319").Value = 15
$ws.Range("R319").Value = "Hortaliza"
